# Auto-generated edit script applying numeric corrections to the
# currentAveragePrice / LevePrice / LeveProfit columns (H:N) across
# several worksheets, per the scheduled runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 159.63637
$ws.Range("I42").Value = 70
$ws.Range("J42").Value = 193.25
$ws.Range("K42").Value = 210
$ws.Range("L42").Value = 579.75
$ws.Range("M42").Value = 20
$ws.Range("N42").Value = -1039.75
$ws.Range("H62").Value = 8463.799999999999
$ws.Range("I62").Value = 7661.8887
$ws.Range("J62").Value = 9666.666999999999
$ws.Range("K62").Value = 7661.8887
$ws.Range("L62").Value = 9666.666999999999
$ws.Range("M62").Value = -7037.8887
$ws.Range("N62").Value = -10914.667
$ws.Range("H65").Value = 8463.799999999999
$ws.Range("I65").Value = 7661.8887
$ws.Range("J65").Value = 9666.666999999999
$ws.Range("K65").Value = 38309.4435
$ws.Range("L65").Value = 48333.335
$ws.Range("M65").Value = -35189.4435
$ws.Range("N65").Value = -54573.335
$ws.Range("H125").Value = 2131.6
$ws.Range("I125").Value = 1470.5714
$ws.Range("K125").Value = 13235.1426
$ws.Range("M125").Value = -10775.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13051.175
$ws.Range("I32").Value = 2924.6155
$ws.Range("J32").Value = 29506.834
$ws.Range("K32").Value = 2924.6155
$ws.Range("L32").Value = 29506.834
$ws.Range("M32").Value = -2637.6155
$ws.Range("N32").Value = -30080.834
$ws.Range("H61").Value = 2059.5908
$ws.Range("I61").Value = 1874.875
$ws.Range("J61").Value = 2552.1667
$ws.Range("K61").Value = 1874.875
$ws.Range("L61").Value = 2552.1667
$ws.Range("M61").Value = -1662.875
$ws.Range("N61").Value = -2976.1667
$ws.Range("H74").Value = 10632910
$ws.Range("I74").Value = 10227678
$ws.Range("J74").Value = 11906496
$ws.Range("K74").Value = 10227678
$ws.Range("L74").Value = 11906496
$ws.Range("M74").Value = -10226804
$ws.Range("N74").Value = -11908244
$ws.Range("H77").Value = 10632910
$ws.Range("I77").Value = 10227678
$ws.Range("J77").Value = 11906496
$ws.Range("K77").Value = 51138390
$ws.Range("L77").Value = 59532480
$ws.Range("M77").Value = -51134022
$ws.Range("N77").Value = -59541216
$ws.Range("H102").Value = 2405.1667
$ws.Range("I102").Value = 2525.4
$ws.Range("J102").Value = 2164.7
$ws.Range("K102").Value = 2525.4
$ws.Range("L102").Value = 2164.7
$ws.Range("M102").Value = -903.4000000000001
$ws.Range("N102").Value = -5408.7
$ws.Range("H132").Value = 2252
$ws.Range("I132").Value = 1608.4706
$ws.Range("J132").Value = 3246.5454
$ws.Range("K132").Value = 4825.4118
$ws.Range("L132").Value = 9739.636200000001
$ws.Range("M132").Value = -2295.4118
$ws.Range("N132").Value = -14799.6362
$ws.Range("H136").Value = 2059.5908
$ws.Range("I136").Value = 1874.875
$ws.Range("J136").Value = 2552.1667
$ws.Range("K136").Value = 5624.625
$ws.Range("L136").Value = 7656.500100000001
$ws.Range("M136").Value = -3074.625
$ws.Range("N136").Value = -12756.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 875.5
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 1501
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 1501
$ws.Range("M22").Value = -77
$ws.Range("N22").Value = -1847
$ws.Range("H99").Value = 1160.0714
$ws.Range("I99").Value = 1214
$ws.Range("J99").Value = 1025.25
$ws.Range("K99").Value = 1214
$ws.Range("L99").Value = 1025.25
$ws.Range("M99").Value = 284
$ws.Range("N99").Value = -4021.25
$ws.Range("H122").Value = 26640
$ws.Range("J122").Value = 26640
$ws.Range("L122").Value = 26640
$ws.Range("N122").Value = -36440
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 31800
$ws.Range("J124").Value = 31800
$ws.Range("L124").Value = 31800
$ws.Range("N124").Value = -41620
$ws.Range("H125").Value = 39999
$ws.Range("J125").Value = 39999
$ws.Range("L125").Value = 39999
$ws.Range("N125").Value = -49839
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2220.8572
$ws.Range("I132").Value = 1145.4546
$ws.Range("J132").Value = 3403.8
$ws.Range("K132").Value = 3436.3638
$ws.Range("L132").Value = 10211.4
$ws.Range("M132").Value = -906.3638000000001
$ws.Range("N132").Value = -15271.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 7964
$ws.Range("I63").Value = 4649.6
$ws.Range("J63").Value = 16250
$ws.Range("K63").Value = 13948.8
$ws.Range("L63").Value = 48750
$ws.Range("M63").Value = -13199.8
$ws.Range("N63").Value = -50248
$ws.Range("H64").Value = 4752.227
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 4835.6665
$ws.Range("K64").Value = 9000
$ws.Range("L64").Value = 14506.9995
$ws.Range("M64").Value = -8730
$ws.Range("N64").Value = -15046.9995
$ws.Range("H66").Value = 7964
$ws.Range("I66").Value = 4649.6
$ws.Range("J66").Value = 16250
$ws.Range("K66").Value = 41846.4
$ws.Range("L66").Value = 146250
$ws.Range("M66").Value = -38102.4
$ws.Range("N66").Value = -153738
$ws.Range("H67").Value = 4752.227
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 4835.6665
$ws.Range("K67").Value = 9000
$ws.Range("L67").Value = 14506.9995
$ws.Range("M67").Value = -8064
$ws.Range("N67").Value = -16378.9995
$ws.Range("H68").Value = 1136.1111
$ws.Range("I68").Value = 756.8333
$ws.Range("K68").Value = 2270.4999
$ws.Range("M68").Value = -1459.4999
$ws.Range("H71").Value = 1136.1111
$ws.Range("I71").Value = 756.8333
$ws.Range("K71").Value = 6811.4997
$ws.Range("M71").Value = -2755.4997
$ws.Range("H107").Value = 941.2222
$ws.Range("I107").Value = 338.8
$ws.Range("J107").Value = 1172.9231
$ws.Range("K107").Value = 1016.4
$ws.Range("L107").Value = 3518.7693
$ws.Range("M107").Value = 903.5999999999999
$ws.Range("N107").Value = -7358.7693
$ws.Range("H114").Value = 1043.2
$ws.Range("I114").Value = 1074
$ws.Range("J114").Value = 1014.7692
$ws.Range("K114").Value = 3222
$ws.Range("L114").Value = 3044.3076
$ws.Range("M114").Value = 32
$ws.Range("N114").Value = -9552.3076
$ws.Range("H127").Value = 1065.5714
$ws.Range("J127").Value = 1065.5714
$ws.Range("L127").Value = 3196.7142
$ws.Range("N127").Value = -13116.7142
$ws.Range("H131").Value = 1093740.1
$ws.Range("I131").Value = 2899172.8
$ws.Range("J131").Value = 978.2368
$ws.Range("K131").Value = 8697518.399999999
$ws.Range("L131").Value = 2934.7104
$ws.Range("M131").Value = -8692478.399999999
$ws.Range("N131").Value = -13014.7104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2852.9412
$ws.Range("I80").Value = 2850
$ws.Range("J80").Value = 2854.5454
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 2854.5454
$ws.Range("M80").Value = -1852
$ws.Range("N80").Value = -4850.5454
$ws.Range("H83").Value = 2852.9412
$ws.Range("I83").Value = 2850
$ws.Range("J83").Value = 2854.5454
$ws.Range("K83").Value = 14250
$ws.Range("L83").Value = 14272.727
$ws.Range("M83").Value = -9258
$ws.Range("N83").Value = -24256.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3004.913
$ws.Range("I7").Value = 2851
$ws.Range("J7").Value = 3559
$ws.Range("K7").Value = 2851
$ws.Range("L7").Value = 3559
$ws.Range("M7").Value = -2739
$ws.Range("N7").Value = -3783
$ws.Range("H126").Value = 3004.913
$ws.Range("I126").Value = 2851
$ws.Range("J126").Value = 3559
$ws.Range("K126").Value = 8553
$ws.Range("L126").Value = 10677
$ws.Range("M126").Value = -6083
$ws.Range("N126").Value = -15617

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 26014.2
$ws.Range("J69").Value = 26014.2
$ws.Range("L69").Value = 26014.2
$ws.Range("N69").Value = -27512.2
$ws.Range("H72").Value = 26014.2
$ws.Range("J72").Value = 26014.2
$ws.Range("L72").Value = 78042.60000000001
$ws.Range("N72").Value = -85530.60000000001
